# Automatische test-sync: 2025-06-27 22:44:50
#
# Appends the new test-mail row (#2, "Kun je 10 dozen schroeven bestellen?")
# to the "Logs" sheet, bumps the matching "Bestelling / Levering" tally on
# the "Dashboard" sheet, and extends the conditional formatting + bar-chart
# series ranges so they keep covering the newly added data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 11
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A11").Value = "Kun je 10 dozen schroeven bestellen?"
$logs.Range("B11").Value = "mailmind.test@zohomail.eu"
$logs.Range("C11").Value = "Testmail #2: Kun je 10 dozen schroeven bestellen?"
$logs.Range("D11").Value = "Bestelling / Levering"

$antwoord11 = @"
Geachte klant,
Bedankt voor uw e-mail. Helaas kan ik geen bestellingen plaatsen, maar ik kan u doorverwijzen naar het bestelteam binnen ons bedrijf. Graag ontvang ik de contactgegevens van uw bedrijf, zodat ik de juiste persoon met u in contact kan brengen.
Ik zie uw reactie graag tegemoet.
Met vriendelijke groet,
[Naam]
E-mailassistent
"@
$logs.Range("E11").Value = $antwoord11

$logs.Range("F11").Value = "2025-06-27 22:44:46"
$logs.Range("G11").Value = "Ja"
$logs.Range("H11").Value = "Ja"
$logs.Range("I11").Value = "Nee"

# Re-flatten the row height: assigning the multi-line E11 text auto-grows
# the row, but row 11 should stay at the sheet's normal (non-custom)
# height, same as every other data row.
$logs.Rows.Item(11).EntireRow.AutoFit()

# ---------------------------------------------------------------------
# 2. Logs sheet: extend the conditional-formatting ranges from row 10 to
#    row 11 (one ModifyAppliesToRange per rule-group is enough since all
#    cfRules sharing a <conditionalFormatting> block move together)
# ---------------------------------------------------------------------
$logs.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D11"))
$logs.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G11"))
$logs.Range("H2:H10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H11"))
$logs.Range("I2:I10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I11"))

# ---------------------------------------------------------------------
# 3. Dashboard sheet: append the tally row for "Bestelling / Levering"
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("A5").Value = "Bestelling / Levering"
$dash.Range("B5").Value = 1

# ---------------------------------------------------------------------
# 4. Dashboard sheet: extend the bar chart's category/value series so it
#    keeps plotting through the newly added row 5
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart()
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$5,'Dashboard'!`$B`$2:`$B`$5,1)"
